$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (title cell A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 00:55"

# Re-order country rows (data re-sorted upstream): swap country labels
$ws.Range("A131").Value = "Siria"
$ws.Range("A132").Value = "Hong Kong"

$ws.Range("A147").Value = "Guyana"
$ws.Range("A148").Value = "Somalia"

$ws.Range("A216").Value = "Montserrat"
$ws.Range("A217").Value = "Islas Malvinas"

# Update numeric statistics per row (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 4
$ws.Range("B4").Value = 8737161
$ws.Range("C4").Value = 71418
$ws.Range("D4").Value = 5692793
$ws.Range("E4").Value = 2815151
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 836
$ws.Range("H4").Value = 229217

# Row 5
$ws.Range("B5").Value = 7813668
$ws.Range("C5").Value = 54028
$ws.Range("D5").Value = 7013569
$ws.Range("E5").Value = 682107
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 656
$ws.Range("H5").Value = 117992

# Row 6
$ws.Range("B6").Value = 5353656
$ws.Range("C6").Value = 21022
$ws.Range("D6").Value = 4797872
$ws.Range("E6").Value = 399313
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 509
$ws.Range("H6").Value = 156471

# Row 9
$ws.Range("B9").Value = 1069368
$ws.Range("C9").Value = 15718
$ws.Range("D9").Value = 866695
$ws.Range("E9").Value = 174335
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 381
$ws.Range("H9").Value = 28338

# Row 11
$ws.Range("B11").Value = 998942
$ws.Range("C11").Value = 8672
$ws.Range("D11").Value = 901652
$ws.Range("E11").Value = 67488
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 166
$ws.Range("H11").Value = 29802

# Row 12
$ws.Range("B12").Value = 883116
$ws.Range("C12").Value = 3240
$ws.Range("D12").Value = 800480
$ws.Range("E12").Value = 48603
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 49
$ws.Range("H12").Value = 34033

# Row 20
$ws.Range("B20").Value = 417350
$ws.Range("C20").Value = 13476
$ws.Range("D20").Value = 310200
$ws.Range("E20").Value = 97060
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 46
$ws.Range("H20").Value = 10090

# Row 33
$ws.Range("B33").Value = 211508
$ws.Range("C33").Value = 2360
$ws.Range("D33").Value = 177746
$ws.Range("E33").Value = 23874
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 26
$ws.Range("H33").Value = 9888

# Row 36
$ws.Range("B36").Value = 158270
$ws.Range("C36").Value = 1819
$ws.Range("D36").Value = 134187
$ws.Range("E36").Value = 11555
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 28
$ws.Range("H36").Value = 12528

# Row 58
$ws.Range("B58").Value = 79574
$ws.Range("C58").Value = 363
$ws.Range("D58").Value = 76143
$ws.Range("E58").Value = 3120
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 311

# Row 62
$ws.Range("B62").Value = 68479
$ws.Range("C62").Value = 1452
$ws.Range("D62").Value = 32412
$ws.Range("E62").Value = 35508
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 7
$ws.Range("H62").Value = 559

# Row 64
$ws.Range("B64").Value = 61882
$ws.Range("C64").Value = 77
$ws.Range("D64").Value = 57190
$ws.Range("E64").Value = 3563
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 1129

# Row 84
$ws.Range("B84").Value = 36519
$ws.Range("C84").Value = 1589
$ws.Range("D84").Value = 18102
$ws.Range("E84").Value = 17340
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 13
$ws.Range("H84").Value = 1077

# Row 112
$ws.Range("B112").Value = 10351
$ws.Range("C112").Value = 9
$ws.Range("D112").Value = 9995
$ws.Range("E112").Value = 287
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 69

# Row 114
$ws.Range("B114").Value = 9015
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = 7361
$ws.Range("E114").Value = 1423
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 231

# Row 131
$ws.Range("B131").Value = 5319
$ws.Range("C131").Value = 52
$ws.Range("D131").Value = 1692
$ws.Range("E131").Value = 3363
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 4
$ws.Range("H131").Value = 264

# Row 132
$ws.Range("B132").Value = 5285
$ws.Range("C132").Value = 4
$ws.Range("D132").Value = 5029
$ws.Range("E132").Value = 151
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 105

# Row 141
$ws.Range("B141").Value = 4401
$ws.Range("C141").Value = 12
$ws.Range("D141").Value = 4160
$ws.Range("E141").Value = 205
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 36

# Row 147
$ws.Range("B147").Value = 3960
$ws.Range("C147").Value = 83
$ws.Range("D147").Value = 2923
$ws.Range("E147").Value = 920
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 117

# Row 148
$ws.Range("B148").Value = 3897
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 3166
$ws.Range("E148").Value = 629
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 102

# Row 161
$ws.Range("B161").Value = 2162
$ws.Range("C161").Value = 23
$ws.Range("D161").Value = 1586
$ws.Range("E161").Value = 524
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 52

# Row 171
$ws.Range("B171").Value = 818
$ws.Range("C171").Value = 14
$ws.Range("D171").Value = 534
$ws.Range("E171").Value = 283
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 1

# Row 172
$ws.Range("B172").Value = 776
$ws.Range("C172").Value = 7
$ws.Range("D172").Value = 702
$ws.Range("E172").Value = 52
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 22

# Row 174
$ws.Range("B174").Value = 699
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 689
$ws.Range("E174").Value = 4
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 6

# Row 191
$ws.Range("B191").Value = 239
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 215
$ws.Range("E191").Value = 23
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 1

# Row 216
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1

# Row 217
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 13
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 0
